# Add viewing of the existing calibration model's plate-reader metadata
# alongside the GFP concentration/RFU calibration data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A for the fluorescent-protein label, and two new
#     trailing columns (D, E) for plate-reader name and gain. ---
$ws.Columns("A").Insert()

# Rename the (now shifted) headers: old A ("Concentration (uM)") is now B,
# old B ("RFU") is now C.
$ws.Range("B1").Value = "Concentration_uM"
$ws.Range("C1").Value = "RFU"

# New column A: protein type
$ws.Range("A1").Value = "FP"
$ws.Range("A2:A7").Value = "GFP"

# New columns D/E: plate reader + gain metadata
$ws.Range("D1").Value = "Plate_Reader"
$ws.Range("D2:D7").Value = "Synergy H1_22060313"

$ws.Range("E1").Value = "Gain"
$ws.Range("E2:E7").Value = 50

# Update selection to match the authored workbook state.
$ws.Range("H9").Select() | Out-Null
